$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 447.5
$ws.Range("I4").Value = 368.57144
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 368.57144
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -254.57144
$ws.Range("N4").Value = -1228

$ws.Range("H12").Value = 3401560.2
$ws.Range("I12").Value = 10204081
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 10204081
$ws.Range("L12").Value = 300
$ws.Range("M12").Value = -10203911
$ws.Range("N12").Value = -640

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1040.909
$ws.Range("I2").Value = 858.7143
$ws.Range("J2").Value = 1359.75
$ws.Range("K2").Value = 858.7143
$ws.Range("L2").Value = 1359.75
$ws.Range("M2").Value = -745.7143
$ws.Range("N2").Value = -1585.75

$ws.Range("H97").Value = 1430.8276
$ws.Range("I97").Value = 1729.6957
$ws.Range("J97").Value = 285.16666
$ws.Range("K97").Value = 1729.6957
$ws.Range("L97").Value = 285.16666
$ws.Range("M97").Value = -1233.6957
$ws.Range("N97").Value = -1277.16666

$ws.Range("H101").Value = 21444.445
$ws.Range("J101").Value = 21444.445
$ws.Range("L101").Value = 21444.445
$ws.Range("N101").Value = -27934.445

$ws.Range("H110").Value = 17926.5
$ws.Range("I110").Value = 19773.143
$ws.Range("K110").Value = 19773.143
$ws.Range("M110").Value = -17728.143

$ws.Range("H116").Value = 1040.909
$ws.Range("I116").Value = 858.7143
$ws.Range("J116").Value = 1359.75
$ws.Range("K116").Value = 858.7143
$ws.Range("L116").Value = 1359.75
$ws.Range("M116").Value = 1435.2857
$ws.Range("N116").Value = -5947.75

$ws.Range("H132").Value = 2099.6597
$ws.Range("I132").Value = 1829.7894
$ws.Range("K132").Value = 5489.3682
$ws.Range("M132").Value = -2959.3682

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1040.909
$ws.Range("I3").Value = 858.7143
$ws.Range("J3").Value = 1359.75
$ws.Range("K3").Value = 858.7143
$ws.Range("L3").Value = 1359.75
$ws.Range("M3").Value = -744.7143
$ws.Range("N3").Value = -1587.75

$ws.Range("H86").Value = 2227.9412
$ws.Range("I86").Value = 2412.5454
$ws.Range("J86").Value = 1889.5
$ws.Range("K86").Value = 2412.5454
$ws.Range("L86").Value = 1889.5
$ws.Range("M86").Value = -1289.5454
$ws.Range("N86").Value = -4135.5

$ws.Range("H89").Value = 2227.9412
$ws.Range("I89").Value = 2412.5454
$ws.Range("J89").Value = 1889.5
$ws.Range("K89").Value = 12062.727
$ws.Range("L89").Value = 9447.5
$ws.Range("M89").Value = -6446.726999999999
$ws.Range("N89").Value = -20679.5

$ws.Range("H94").Value = 13201.8125
$ws.Range("I94").Value = 507.9091
$ws.Range("J94").Value = 41128.4
$ws.Range("K94").Value = 507.9091
$ws.Range("L94").Value = 41128.4
$ws.Range("M94").Value = -56.90910000000002
$ws.Range("N94").Value = -42030.4

$ws.Range("H105").Value = 3621.6155
$ws.Range("I105").Value = 4301.8
$ws.Range("J105").Value = 3196.5
$ws.Range("K105").Value = 4301.8
$ws.Range("L105").Value = 3196.5
$ws.Range("M105").Value = -2554.8
$ws.Range("N105").Value = -6690.5

$ws.Range("H107").Value = 1462.6666
$ws.Range("I107").Value = 1216.2
$ws.Range("K107").Value = 1216.2
$ws.Range("M107").Value = 703.8

$ws.Range("H134").Value = 2232.3635
$ws.Range("I134").Value = 1795.742
$ws.Range("K134").Value = 5387.226
$ws.Range("M134").Value = -2852.226

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 1419.75
$ws.Range("I35").Value = 675
$ws.Range("J35").Value = 3654
$ws.Range("K35").Value = 675
$ws.Range("L35").Value = 3654
$ws.Range("M35").Value = -381
$ws.Range("N35").Value = -4242

$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()

$ws.Range("H62").Value = 90913990
$ws.Range("I62").Value = 3140
$ws.Range("J62").Value = 166673020
$ws.Range("K62").Value = 3140
$ws.Range("L62").Value = 166673020
$ws.Range("M62").Value = -2516
$ws.Range("N62").Value = -166674268

$ws.Range("H65").Value = 90913990
$ws.Range("I65").Value = 3140
$ws.Range("J65").Value = 166673020
$ws.Range("K65").Value = 15700
$ws.Range("L65").Value = 833365100
$ws.Range("M65").Value = -12580
$ws.Range("N65").Value = -833371340

$ws.Range("H105").Value = 1557.6666
$ws.Range("I105").Value = 1708
$ws.Range("J105").Value = 1181.8334
$ws.Range("K105").Value = 1708
$ws.Range("L105").Value = 1181.8334
$ws.Range("M105").Value = 39
$ws.Range("N105").Value = -4675.8334

$ws.Range("H129").Value = 32499.5
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 32499.5
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 32499.5
$ws.Range("N129").Value = -42499.5
$ws.Range("M129").ClearContents()

$ws.Range("H132").Value = 3770.7144
$ws.Range("I132").Value = 2861.5
$ws.Range("J132").Value = 4983
$ws.Range("K132").Value = 8584.5
$ws.Range("L132").Value = 14949
$ws.Range("M132").Value = -6054.5
$ws.Range("N132").Value = -20009

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1281.8572
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 1328.8334
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 3986.5002
$ws.Range("M80").Value = -2064
$ws.Range("N80").Value = -5858.5002

$ws.Range("H83").Value = 1281.8572
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 1328.8334
$ws.Range("K83").Value = 9000
$ws.Range("L83").Value = 11959.5006
$ws.Range("M83").Value = -4320
$ws.Range("N83").Value = -21319.5006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 22.153847
$ws.Range("I2").Value = 15.777778
$ws.Range("J2").Value = 36.5
$ws.Range("K2").Value = 15.777778
$ws.Range("L2").Value = 36.5
$ws.Range("M2").Value = 97.222222
$ws.Range("N2").Value = -262.5

$ws.Range("H97").Value = 2002.5
$ws.Range("I97").Value = 2002.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2002.5
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1506.5
$ws.Range("N97").ClearContents()

$ws.Range("H113").Value = 2813.5833
$ws.Range("I113").Value = 1837.75
$ws.Range("J113").Value = 3301.5
$ws.Range("K113").Value = 1837.75
$ws.Range("L113").Value = 3301.5
$ws.Range("M113").Value = 332.25
$ws.Range("N113").Value = -7641.5

$ws.Range("H122").Value = 1736.2222
$ws.Range("I122").Value = 1683.4667
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5050.4001
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2600.4001
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 2520.0527
$ws.Range("I132").Value = 2044.8462
$ws.Range("J132").Value = 3549.6667
$ws.Range("K132").Value = 6134.5386
$ws.Range("L132").Value = 10649.0001
$ws.Range("M132").Value = -3604.5386
$ws.Range("N132").Value = -15709.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 799.3077
$ws.Range("I22").Value = 498.5
$ws.Range("J22").Value = 1057.1428
$ws.Range("K22").Value = 498.5
$ws.Range("L22").Value = 1057.1428
$ws.Range("M22").Value = -203.5
$ws.Range("N22").Value = -1647.1428

$ws.Range("H27").Value = 799.3077
$ws.Range("I27").Value = 498.5
$ws.Range("J27").Value = 1057.1428
$ws.Range("K27").Value = 498.5
$ws.Range("L27").Value = 1057.1428
$ws.Range("M27").Value = -391.5
$ws.Range("N27").Value = -1271.1428

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 17779
$ws.Range("J64").Value = 17779
$ws.Range("L64").Value = 17779
$ws.Range("N64").Value = -18275

$ws.Range("H67").Value = 17779
$ws.Range("J67").Value = 17779
$ws.Range("L67").Value = 17779
$ws.Range("N67").Value = -19495

$ws.Range("H122").Value = 47621052
$ws.Range("I122").Value = 100001400
$ws.Range("J122").Value = 2554.5454
$ws.Range("K122").Value = 300004200
$ws.Range("L122").Value = 7663.6362
$ws.Range("M122").Value = -300001750
$ws.Range("N122").Value = -12563.6362

$ws.Range("H136").Value = 1459.7742
$ws.Range("I136").Value = 1171.4348
$ws.Range("J136").Value = 2288.75
$ws.Range("K136").Value = 3514.3044
$ws.Range("L136").Value = 6866.25
$ws.Range("M136").Value = -964.3044
$ws.Range("N136").Value = -11966.25
